$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 53 with the new titration data
$ws.Range("A53").Value = 20220127
$ws.Range("B53").Value = 2231.4175319214601
$ws.Range("C53").Value = 2224.4699999999998
$ws.Range("D53").Formula = "=100*(B53-C53)/C53"
$ws.Range("E53").Value = 180
$ws.Range("F53").Value = "CRM OPENED 20220118"

# Update selection to reflect new active cell after data entry
$ws.Range("F54").Select()
